# Update "想去人数" (want-to-go count) column F values on both the
# "展览" and "全部类型" worksheets, rows 2-9.

$wb = $excel.ActiveWorkbook

$newValues = @{
    2 = 2293
    3 = 1744
    4 = 339
    5 = 1097
    6 = 883
    7 = 42
    8 = 5856
    9 = 89
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $newValues.Keys) {
        $ws.Cells.Item($row, 6).Value = $newValues[$row]
    }
}
